# Fix typo on slide 7 ("Role" slide): "Yehor Hora" -> "Yehor Hoda"
# (matches the spelling used elsewhere in the deck, e.g. the title slide's
# credits list already reads "... Yehor Hoda").
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$found = $tr.Find("Hora")
if ($found -ne $null) {
    $start = $found.Start
    # "Hora" -> change the third letter ("r") to "d", giving "Hoda"
    $rChar = $tr.Characters($start + 2, 1)
    $rChar.Text = "d"
}
